$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# "Periodo Mora" (column E) values for rows 16-28 get reversed (the previously
# used account-statement periods are dropped, new ones are added, effectively
# reversing the period list top-to-bottom).
$ws.Range("E16").Value = "2112"
$ws.Range("E17").Value = "2111"
$ws.Range("E18").Value = "2110"
$ws.Range("E19").Value = "2109"
$ws.Range("E20").Value = "2108"
$ws.Range("E21").Value = "2107"
$ws.Range("E22").Value = "2106"
$ws.Range("E23").Value = "2105"
$ws.Range("E24").Value = "2104"
$ws.Range("E25").Value = "2103"
$ws.Range("E26").Value = "2102"
$ws.Range("E27").Value = "2101"
$ws.Range("E28").Value = "2012"

# "Valor Mora" (column F) values for the first and last data row swap along
# with the period reordering above.
$ws.Range("F16").Value = 21333
$ws.Range("F28").Value = 40000
